$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 640, shifting existing data (old rows 640-692) down to 644-696
$ws.Rows("640:643").Insert()

# New row 640: Especial
$ws.Cells.Item(640, 1).Value = 2
$ws.Cells.Item(640, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(640, 3).Value = "Coquimbo"
$ws.Cells.Item(640, 4).Value = 45106
$ws.Cells.Item(640, 5).Value = 4
$ws.Cells.Item(640, 6).Value = 100112043
$ws.Cells.Item(640, 7).Value = "Pepino dulce"
$ws.Cells.Item(640, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(640, 9).Value = "Especial"
$ws.Cells.Item(640, 10).Value = 200
$ws.Cells.Item(640, 11).Value = 11000
$ws.Cells.Item(640, 12).Value = 12000
$ws.Cells.Item(640, 13).Value = 11500
$ws.Cells.Item(640, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(640, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(640, 16).Value = 639
$ws.Cells.Item(640, 17).Value = 18
$ws.Cells.Item(640, 18).Value = "Hortaliza"

# New row 641: Primera
$ws.Cells.Item(641, 1).Value = 2
$ws.Cells.Item(641, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(641, 3).Value = "Coquimbo"
$ws.Cells.Item(641, 4).Value = 45106
$ws.Cells.Item(641, 5).Value = 4
$ws.Cells.Item(641, 6).Value = 100112043
$ws.Cells.Item(641, 7).Value = "Pepino dulce"
$ws.Cells.Item(641, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(641, 9).Value = "Primera"
$ws.Cells.Item(641, 10).Value = 400
$ws.Cells.Item(641, 11).Value = 9000
$ws.Cells.Item(641, 12).Value = 10000
$ws.Cells.Item(641, 13).Value = 9500
$ws.Cells.Item(641, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(641, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(641, 16).Value = 528
$ws.Cells.Item(641, 17).Value = 18
$ws.Cells.Item(641, 18).Value = "Hortaliza"

# New row 642: Segunda
$ws.Cells.Item(642, 1).Value = 2
$ws.Cells.Item(642, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(642, 3).Value = "Coquimbo"
$ws.Cells.Item(642, 4).Value = 45106
$ws.Cells.Item(642, 5).Value = 4
$ws.Cells.Item(642, 6).Value = 100112043
$ws.Cells.Item(642, 7).Value = "Pepino dulce"
$ws.Cells.Item(642, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(642, 9).Value = "Segunda"
$ws.Cells.Item(642, 10).Value = 240
$ws.Cells.Item(642, 11).Value = 7000
$ws.Cells.Item(642, 12).Value = 8000
$ws.Cells.Item(642, 13).Value = 7500
$ws.Cells.Item(642, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(642, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(642, 16).Value = 417
$ws.Cells.Item(642, 17).Value = 18
$ws.Cells.Item(642, 18).Value = "Hortaliza"

# New row 643: Tercera
$ws.Cells.Item(643, 1).Value = 2
$ws.Cells.Item(643, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(643, 3).Value = "Coquimbo"
$ws.Cells.Item(643, 4).Value = 45106
$ws.Cells.Item(643, 5).Value = 4
$ws.Cells.Item(643, 6).Value = 100112043
$ws.Cells.Item(643, 7).Value = "Pepino dulce"
$ws.Cells.Item(643, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(643, 9).Value = "Tercera"
$ws.Cells.Item(643, 10).Value = 200
$ws.Cells.Item(643, 11).Value = 4000
$ws.Cells.Item(643, 12).Value = 5000
$ws.Cells.Item(643, 13).Value = 4500
$ws.Cells.Item(643, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(643, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(643, 16).Value = 250
$ws.Cells.Item(643, 17).Value = 18
$ws.Cells.Item(643, 18).Value = "Hortaliza"
